# Commit: "Sanity checks and a whole bunch of upgrades"
# Applies the new submissions (#47-#53) to the "Submissions" sheet, and
# fixes the off-by-one label on the existing #46 entry (was mislabeled
# "(45)+..." -> "(44)+...").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Submissions")
$ws.Activate()

# NOTE: new text written to a cell is appended to the shared-string table in
# the order the writes happen (existing repeated strings, like "No", are
# reused). To reproduce the exact shared-string indices from the original
# commit, the new unique strings below are written in this precise order:
#   D49, D50, D51, D52, Q50, D53, D54, D55

# --- Fix existing row 48 (id 46) and fill in its previously-blank columns ---
$ws.Range("D48").Value = "(44)+tweak loader, fallbacks in simple, unlock t0"

# --- Row 49 (id 47) ---
$ws.Range("D49").Value = "(46)+regularization in simple and GP"

# --- Row 50 (id 48) ---
$ws.Range("D50").Value = "(47)+force Kepler"

# --- Row 51 (id 49) ---
$ws.Range("D51").Value = "(48)+fudge based on var"

# --- Row 52 (id 50) ---
$ws.Range("D52").Value = "(49)+adjust based on u"

# --- Row 50 conclusion note (written after D51/D52 to match string order) ---
$ws.Range("Q50").Value = "Not using force Kepler"

# --- Row 53 (id 51) ---
$ws.Range("D53").Value = "(47)+fudge based on var"

# --- Row 54 (id 52) ---
$ws.Range("D54").Value = "(51)+adjust based on u"

# --- Row 55 (id 53) ---
$ws.Range("D55").Value = "(52)+fudge based on multi, multi-transit"

# --- Remaining columns (numbers / reused strings - order independent) ---
$ws.Range("E48").Value = "No"
$ws.Range("G48").Value = 0.6122
$ws.Range("P48").Value = 0.61
$ws.Range("P48").NumberFormat = "0.000"

$ws.Range("A49").Value = 47
$ws.Range("E49").Value = "No"
$ws.Range("G49").Value = 0.6214
$ws.Range("P49").Value = 0.615
$ws.Range("P49").NumberFormat = "0.000"

$ws.Range("A50").Value = 48
$ws.Range("E50").Value = "No"
$ws.Range("G50").Value = 0.6228
$ws.Range("P50").Value = 0.612
$ws.Range("P50").NumberFormat = "0.000"

$ws.Range("A51").Value = 49
$ws.Range("E51").Value = "No"
$ws.Range("G51").Value = 0.6314
$ws.Range("P51").NumberFormat = "0.000"

$ws.Range("A52").Value = 50
$ws.Range("E52").Value = "No"
$ws.Range("G52").Value = 0.6363
$ws.Range("P52").NumberFormat = "0.000"

$ws.Range("A53").Value = 51
$ws.Range("E53").Value = "No"
$ws.Range("G53").Value = 0.6303
$ws.Range("P53").Value = 0.618
$ws.Range("P53").NumberFormat = "0.000"

$ws.Range("A54").Value = 52
$ws.Range("E54").Value = "No"
$ws.Range("G54").Value = 0.6352
$ws.Range("P54").Value = 0.622
$ws.Range("P54").NumberFormat = "0.000"

$ws.Range("A55").Value = 53
$ws.Range("E55").Value = "No"
$ws.Range("G55").Value = 0.6399
$ws.Range("P55").NumberFormat = "0.000"

# --- Extend Table2 to cover the new rows ---
$wb.Worksheets.Item("Submissions").ListObjects.Item("Table2").Resize($ws.Range("A1:Q55"))

# --- Match the recorded selection on this sheet ---
$ws.Range("Q52").Select()

Write-Output "done"
